# The "codebook" sheet template had a stray, unused formatted cell in L1
# (trailing empty column next to the K1 border-cap cell). Clear it so it
# no longer renders as a separate bordered cell.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("codebook")
$ws1.Range("L1").Clear()

# The "additional information" sheet listed "Number of Indicators" twice
# (rows 8 and 10) - a duplicate variable name. Remove the duplicate row;
# the rows below it (9 and 10) shift up to fill the gap.
$ws3 = $wb.Worksheets.Item("additional information")
$ws3.Rows.Item(8).Delete()

# Restore the cursor/selection state for each sheet, ending with
# "additional information" as the active tab (matching the saved file).
[void]$ws1.Range("C25").Select()
[void]$ws3.Range("F17").Select()
